$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.237.76"
$ws.Range("E2").Value = "  -2.84%  "
$ws.Range("D3").Value = "2.655.35"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.89"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.12"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.03"
$ws.Range("E9").Value = "  +9.22%  "
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.335"
$ws.Range("E11").Value = "  -2.33%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "3.124.42"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").Value = "59.220.67"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.18"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000136"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "2.663.30"
$ws.Range("E17").Value = "  -6.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "340.17"
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.37"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.69"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.413"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").Value = "  -2.06%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.67"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "148.87"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.14"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.898"
$ws.Range("E36").Value = "  -6.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.883"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.74"
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.49"
$ws.Range("E39").Value = "  -5.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.60"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.619"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "276.35"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.88"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0973"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0537"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").Value = "  +2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.81"
$ws.Range("E48").Value = "  -2.77%  "
$ws.Range("D49").Value = "2.029.00"
$ws.Range("E49").Value = "  -5.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.99"
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("E51").Value = "  -2.78%  "
